$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.939310073852539
$ws.Range("B1").Value = 4.401539325714111
$ws.Range("C1").Value = 3.850996971130371
$ws.Range("D1").Value = 4.774041175842285
$ws.Range("E1").Value = 4.836911201477051
